$wb = $excel.ActiveWorkbook

$wsRetained = $wb.Worksheets.Item("samples_retained")
$wsPositive = $wb.Worksheets.Item("positive")
$wsNegative = $wb.Worksheets.Item("negative")
$wsDiscard  = $wb.Worksheets.Item("discard")

# --- samples_retained sheet: start the new jl-corpus row ---
$wsRetained.Range("A19").Value = "jl-corpus"
$wsRetained.Range("B19").Value = "acted"
$wsRetained.Range("C19").Value = 480
$wsRetained.Range("D19").Value = 960
$wsRetained.Range("E19").Value = 240
$wsRetained.Range("F19").Value = "English"
$wsRetained.Range("G19").Formula = "=IF(OR(ISBLANK(C19), ISBLANK(D19),ISBLANK(E19)), """", SUM(C19:E19))"

# --- positive sheet: add pensive, enthusiastic rows ---
$wsPositive.Range("A11").Value = "pensive"
$wsPositive.Range("B11").Value = "en"
$wsPositive.Range("C11").Value = "pen"
$wsPositive.Range("A12").Value = "enthusiastic"
$wsPositive.Range("B12").Value = "en"
$wsPositive.Range("C12").Value = "ent"

# --- negative sheet: add apologetic, anxious, worried rows ---
$wsNegative.Range("A28").Value = "apologetic"
$wsNegative.Range("B28").Value = "en"
$wsNegative.Range("C28").Value = "apo"
$wsNegative.Range("A29").Value = "anxious"
$wsNegative.Range("B29").Value = "en"
$wsNegative.Range("C29").Value = "anx"
$wsNegative.Range("A30").Value = "worried"
$wsNegative.Range("B30").Value = "en"
$wsNegative.Range("C30").Value = "wor"

# --- positive sheet: add excited row ---
$wsPositive.Range("A13").Value = "excited"
$wsPositive.Range("B13").Value = "en"
$wsPositive.Range("C13").Value = "exc"

# --- discard sheet: add encouraging, concerned, assertive rows ---
$wsDiscard.Range("A5").Value = "encouraging"
$wsDiscard.Range("B5").Value = "en"
$wsDiscard.Range("A6").Value = "concerned"
$wsDiscard.Range("B6").Value = "en"
$wsDiscard.Range("A7").Value = "assertive"
$wsDiscard.Range("B7").Value = "en"

# --- samples_retained sheet: finish the jl-corpus row with its note ---
$wsRetained.Range("H19").Value = "New Zealand English; valence labels were provided for the non-primary emotions; semi-natural elicitation (almost spont.); apologetic, anxious, worried negative; excited and happy positive"

# --- selections to match final cursor positions ---
# (select the non-active sheet's target cell first so the workbook's
# active/"tabSelected" sheet ends up back on samples_retained, matching
# the source file where tabSelected stays on sheet1 throughout)
$wsPositive.Range("A14").Select()
$wsRetained.Range("C20").Select()
